# The "{{ seal_quantity }}" merge-field placeholder (table cell) was
# authored with black, theme-"text1" run colour. Per the commit
# ("added extra checker for excel") the cell is being turned into a
# hidden/white-on-white marker: its 4 runs ("{{ ", "seal", "_quantity",
# " }}") need their run-level colour flipped from
#   000000 / themeColor=text1   ->   FFFFFF / themeColor=background1
# The paragraph mark's own rPr colour (held in pPr/rPr) is left as-is -
# only the run Font colours change.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "{{ seal_quantity }}",  # FindText
    $true,                  # MatchCase
    $true,                  # MatchWholeWord
    $false,                 # MatchWildcards
    $false,                 # MatchSoundsLike
    $false,                 # MatchAllWordForms
    $true,                  # Forward
    1,                      # Wrap (wdFindContinue)
    $false,                 # Format
    "",                     # ReplaceWith
    0)                      # Replace (wdReplaceNone) - formatting only, no text change

if (-not $found) {
    throw "Could not find the '{{ seal_quantity }}' placeholder run"
}

# Recolour just the matched run text (the cell's four runs) to the
# theme's background1 colour (white), mirroring how Word itself writes
# w:themeColor="background1" when a theme swatch is applied.
$rng.Font.TextColor.ObjectThemeColor = 12   # wdThemeColorBackground1
